# Insert a new weekly price record as row 33, pushing the existing
# rows 33-39 down to rows 34-40 (dimension grows from A1:R39 to A1:R40).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("33:33").Insert()

$ws.Range("A33").Value2 = 10
$ws.Range("B33").Value2 = "Vega Modelo de Temuco"
$ws.Range("C33").Value2 = "La Araucanía"
$ws.Range("D33").Value2 = 44798
$ws.Range("E33").Value2 = 9
$ws.Range("F33").Value2 = 100112042
$ws.Range("G33").Value2 = "Locoto"
$ws.Range("H33").Value2 = "Sin especificar"
$ws.Range("I33").Value2 = "Primera"
$ws.Range("J33").Value2 = 80
$ws.Range("K33").Value2 = 2700
$ws.Range("L33").Value2 = 2700
$ws.Range("M33").Value2 = 2700
$ws.Range("N33").Value2 = "$/kilo"
$ws.Range("O33").Value2 = "Región de Arica y Parinacota"
$ws.Range("P33").Value2 = 2700
$ws.Range("Q33").Value2 = 1
$ws.Range("R33").Value2 = "Hortaliza"
